$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so
# numeric-looking strings (e.g. "6.70", "0.680") are not coerced into
# numbers and stripped of trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.938.10"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.915.61"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "593.22"
$ws.Range("E5").Value = "  +1.01%  "

# Row 6 - Solana
$ws.Range("D6").Value = "146.11"
$ws.Range("E6").Value = "  -0.28%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.76%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "6.87"
$ws.Range("E9").Value = "  +2.01%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  -0.32%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -1.86%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  +0.43%  "

# Row 13 - Avalanche
$ws.Range("E13").Value = "  +0.13%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.25%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.398.36"
$ws.Range("E15").Value = "  +0.08%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "60.904.27"
$ws.Range("E16").Value = "  +0.13%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "6.70"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.917.38"
$ws.Range("E18").Value = "  +0.23%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "430.77"
$ws.Range("E19").Value = "  +0.89%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "13.38"
$ws.Range("E20").Value = "  -1.43%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "0.680"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.81%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "81.60"
$ws.Range("E23").Value = "  +1.92%  "

# Row 24 - RenderToken
$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -0.34%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "11.92"
$ws.Range("E26").Value = "  +0.48%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.04%  "

# Row 28 - ImmutableX
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +5.13%  "

# Row 29 - FirstDigitalUSD
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.33%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  -2.90%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "26.49"
$ws.Range("E32").Value = "  +0.15%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  +1.52%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0853"
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +0.30%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  +0.10%  "

# Row 37 - dogwifhat
$ws.Range("D37").Value = "3.03"
$ws.Range("E37").Value = "  +1.32%  "

# Row 38 - was Stacks, now Kaspa
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.123"
$ws.Range("E38").Value = "  -0.02%  "

# Row 39 - was Kaspa, now Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.99"
$ws.Range("E39").Value = "  -1.58%  "

# Row 40 - Cosmos
$ws.Range("D40").Value = "8.56"

# Row 41 - TheGraph
$ws.Range("D41").Value = "0.286"
$ws.Range("E41").Value = "  -2.15%  "

# Row 42 - Arweave
$ws.Range("D42").Value = "39.96"
$ws.Range("E42").Value = "  -4.75%  "

# Row 43 - Bittensor
$ws.Range("D43").Value = "375.65"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0345"
$ws.Range("E44").Value = "  -0.86%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.705.63"
$ws.Range("E45").Value = "  +1.24%  "

# Row 46 - Monero
$ws.Range("D46").Value = "132.11"
$ws.Range("E46").Value = "  -0.66%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "23.89"
$ws.Range("E48").Value = "  -4.42%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.13%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -3.61%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +1.61%  "
